$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.519.00'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '2.894.81'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.65'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.45'
$ws.Range("E6").Value = '  -4.45%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.554'
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("D9").Value = '2.899.87'
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.362'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '3.411.06'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '60.541.58'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.63'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '2.898.90'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000141'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.65'
$ws.Range("E20").Value = '  -0.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '363.46'
$ws.Range("E21").Value = '  -2.35%  '
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.51'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.454'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.181'
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.78'
$ws.Range("E28").Value = '  -5.04%  '
$ws.Range("E29").Value = '  -4.11%  '
$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.68'
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.63'
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.02'
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.37'
$ws.Range("E34").Value = '  -2.65%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.56'
$ws.Range("E35").Value = '  -5.17%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("E36").Value = '  -4.60%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.20'
$ws.Range("E37").Value = '  -3.71%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.47'
$ws.Range("E38").Value = '  +3.31%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.69'
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("D41").Value = '2.280.20'
$ws.Range("E41").Value = '  -4.42%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.645'
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0580'
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.66'
$ws.Range("E44").Value = '  -4.10%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.95'
$ws.Range("E46").Value = '  +3.33%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0236'
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0928'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.32'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '249.72'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.40'
$ws.Range("E51").Value = '  -4.39%  '
